$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lama1"
$ws.Cells.Item(2, 3).Value = "Rpsa"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.014112
$ws.Cells.Item(2, 8).Value = 0.042336
$ws.Cells.Item(2, 9).Value = 0.1773673913134555
$ws.Cells.Item(2, 10).Value = 0.1773673913134555
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 156.53184
$ws.Cells.Item(2, 14).Value = 469.59552
$ws.Cells.Item(2, 15).Value = 0.1403721039197297
$ws.Cells.Item(2, 16).Value = 0.1403721039197297
$ws.Cells.Item(2, 17).Value = 2.20897732608
$ws.Cells.Item(2, 18).Value = 19.88079593472
$ws.Cells.Item(2, 19).Value = 0.02489743388542374
$ws.Cells.Item(2, 20).Value = 0.02489743388542373

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lama1"
$ws.Cells.Item(3, 3).Value = "Rpsa"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.014112
$ws.Cells.Item(3, 8).Value = 0.042336
$ws.Cells.Item(3, 9).Value = 0.1773673913134555
$ws.Cells.Item(3, 10).Value = 0.1773673913134555
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 389.2008666666666
$ws.Cells.Item(3, 14).Value = 1167.6026
$ws.Cells.Item(3, 15).Value = 0.3490212885850074
$ws.Cells.Item(3, 16).Value = 0.3490212885850074
$ws.Cells.Item(3, 17).Value = 5.492402630399999
$ws.Cells.Item(3, 18).Value = 49.43162367359999
$ws.Cells.Item(3, 19).Value = 0.06190499546918348
$ws.Cells.Item(3, 20).Value = 0.06190499546918347

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lama1"
$ws.Cells.Item(4, 3).Value = "Rpsa"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.014112
$ws.Cells.Item(4, 8).Value = 0.042336
$ws.Cells.Item(4, 9).Value = 0.1773673913134555
$ws.Cells.Item(4, 10).Value = 0.1773673913134555
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 169.4499613333333
$ws.Cells.Item(4, 14).Value = 508.349884
$ws.Cells.Item(4, 15).Value = 0.1519566088373896
$ws.Cells.Item(4, 16).Value = 0.1519566088373896
$ws.Cells.Item(4, 17).Value = 2.391277854336
$ws.Cells.Item(4, 18).Value = 21.521500689024
$ws.Cells.Item(4, 19).Value = 0.02695214730232696
$ws.Cells.Item(4, 20).Value = 0.02695214730232696

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Lama1"
$ws.Cells.Item(5, 3).Value = "Rpsa"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.014112
$ws.Cells.Item(5, 8).Value = 0.042336
$ws.Cells.Item(5, 9).Value = 0.1773673913134555
$ws.Cells.Item(5, 10).Value = 0.1773673913134555
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 169.0002543333333
$ws.Cells.Item(5, 14).Value = 507.000763
$ws.Cells.Item(5, 15).Value = 0.1515533278324679
$ws.Cells.Item(5, 16).Value = 0.1515533278324679
$ws.Cells.Item(5, 17).Value = 2.384931589152
$ws.Cells.Item(5, 18).Value = 21.464384302368
$ws.Cells.Item(5, 19).Value = 0.02688061840251775
$ws.Cells.Item(5, 20).Value = 0.02688061840251774

# Row 6
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Lama1"
$ws.Cells.Item(6, 3).Value = "Rpsa"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.014112
$ws.Cells.Item(6, 8).Value = 0.042336
$ws.Cells.Item(6, 9).Value = 0.1773673913134555
$ws.Cells.Item(6, 10).Value = 0.1773673913134555
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 60.30985666666667
$ws.Cells.Item(6, 14).Value = 180.92957
$ws.Cells.Item(6, 15).Value = 0.05408370250677011
$ws.Cells.Item(6, 16).Value = 0.05408370250677011
$ws.Cells.Item(6, 17).Value = 0.85109269728
$ws.Cells.Item(6, 18).Value = 7.659834275520001
$ws.Cells.Item(6, 19).Value = 0.009592685226198808
$ws.Cells.Item(6, 20).Value = 0.009592685226198806

# Row 7
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Lama1"
$ws.Cells.Item(7, 3).Value = "Rpsa"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.014112
$ws.Cells.Item(7, 8).Value = 0.042336
$ws.Cells.Item(7, 9).Value = 0.1773673913134555
$ws.Cells.Item(7, 10).Value = 0.1773673913134555
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 170.6279296666667
$ws.Cells.Item(7, 14).Value = 511.883789
$ws.Cells.Item(7, 15).Value = 0.1530129683186351
$ws.Cells.Item(7, 16).Value = 0.1530129683186351
$ws.Cells.Item(7, 17).Value = 2.407901343456
$ws.Cells.Item(7, 18).Value = 21.671112091104
$ws.Cells.Item(7, 19).Value = 0.02713951102780473
$ws.Cells.Item(7, 20).Value = 0.02713951102780473

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Lama1"
$ws.Cells.Item(8, 3).Value = "Rpsa"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.059263
$ws.Cells.Item(8, 8).Value = 0.177789
$ws.Cells.Item(8, 9).Value = 0.7448500362393221
$ws.Cells.Item(8, 10).Value = 0.7448500362393219
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 156.53184
$ws.Cells.Item(8, 14).Value = 469.59552
$ws.Cells.Item(8, 15).Value = 0.1403721039197297
$ws.Cells.Item(8, 16).Value = 0.1403721039197297
$ws.Cells.Item(8, 17).Value = 9.27654643392
$ws.Cells.Item(8, 18).Value = 83.48891790527999
$ws.Cells.Item(8, 19).Value = 0.1045561666916006
$ws.Cells.Item(8, 20).Value = 0.1045561666916005

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Lama1"
$ws.Cells.Item(9, 3).Value = "Rpsa"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.059263
$ws.Cells.Item(9, 8).Value = 0.177789
$ws.Cells.Item(9, 9).Value = 0.7448500362393221
$ws.Cells.Item(9, 10).Value = 0.7448500362393219
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 389.2008666666666
$ws.Cells.Item(9, 14).Value = 1167.6026
$ws.Cells.Item(9, 15).Value = 0.3490212885850074
$ws.Cells.Item(9, 16).Value = 0.3490212885850074
$ws.Cells.Item(9, 17).Value = 23.06521096126666
$ws.Cells.Item(9, 18).Value = 207.5868986514
$ws.Cells.Item(9, 19).Value = 0.2599685194508377
$ws.Cells.Item(9, 20).Value = 0.2599685194508375

# Row 10
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Lama1"
$ws.Cells.Item(10, 3).Value = "Rpsa"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.059263
$ws.Cells.Item(10, 8).Value = 0.177789
$ws.Cells.Item(10, 9).Value = 0.7448500362393221
$ws.Cells.Item(10, 10).Value = 0.7448500362393219
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 169.4499613333333
$ws.Cells.Item(10, 14).Value = 508.349884
$ws.Cells.Item(10, 15).Value = 0.1519566088373896
$ws.Cells.Item(10, 16).Value = 0.1519566088373896
$ws.Cells.Item(10, 17).Value = 10.04211305849733
$ws.Cells.Item(10, 18).Value = 90.37901752647599
$ws.Cells.Item(10, 19).Value = 0.1131848855993341
$ws.Cells.Item(10, 20).Value = 0.1131848855993341

# Row 11
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Lama1"
$ws.Cells.Item(11, 3).Value = "Rpsa"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.059263
$ws.Cells.Item(11, 8).Value = 0.177789
$ws.Cells.Item(11, 9).Value = 0.7448500362393221
$ws.Cells.Item(11, 10).Value = 0.7448500362393219
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 169.0002543333333
$ws.Cells.Item(11, 14).Value = 507.000763
$ws.Cells.Item(11, 15).Value = 0.1515533278324679
$ws.Cells.Item(11, 16).Value = 0.1515533278324679
$ws.Cells.Item(11, 17).Value = 10.01546207255633
$ws.Cells.Item(11, 18).Value = 90.139158653007
$ws.Cells.Item(11, 19).Value = 0.1128845017282036
$ws.Cells.Item(11, 20).Value = 0.1128845017282036

# Row 12
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Lama1"
$ws.Cells.Item(12, 3).Value = "Rpsa"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.059263
$ws.Cells.Item(12, 8).Value = 0.177789
$ws.Cells.Item(12, 9).Value = 0.7448500362393221
$ws.Cells.Item(12, 10).Value = 0.7448500362393219
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 60.30985666666667
$ws.Cells.Item(12, 14).Value = 180.92957
$ws.Cells.Item(12, 15).Value = 0.05408370250677011
$ws.Cells.Item(12, 16).Value = 0.05408370250677011
$ws.Cells.Item(12, 17).Value = 3.574143035636667
$ws.Cells.Item(12, 18).Value = 32.16728732073
$ws.Cells.Item(12, 19).Value = 0.04028424777212443
$ws.Cells.Item(12, 20).Value = 0.04028424777212442

# Row 13
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Lama1"
$ws.Cells.Item(13, 3).Value = "Rpsa"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.059263
$ws.Cells.Item(13, 8).Value = 0.177789
$ws.Cells.Item(13, 9).Value = 0.7448500362393221
$ws.Cells.Item(13, 10).Value = 0.7448500362393219
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 170.6279296666667
$ws.Cells.Item(13, 14).Value = 511.883789
$ws.Cells.Item(13, 15).Value = 0.1530129683186351
$ws.Cells.Item(13, 16).Value = 0.1530129683186351
$ws.Cells.Item(13, 17).Value = 10.11192299583567
$ws.Cells.Item(13, 18).Value = 91.007306962521
$ws.Cells.Item(13, 19).Value = 0.1139717149972216
$ws.Cells.Item(13, 20).Value = 0.1139717149972216

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Lama1"
$ws.Cells.Item(14, 3).Value = "Rpsa"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.006188666666666666
$ws.Cells.Item(14, 8).Value = 0.018566
$ws.Cells.Item(14, 9).Value = 0.07778257244722256
$ws.Cells.Item(14, 10).Value = 0.07778257244722254
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 156.53184
$ws.Cells.Item(14, 14).Value = 469.59552
$ws.Cells.Item(14, 15).Value = 0.1403721039197297
$ws.Cells.Item(14, 16).Value = 0.1403721039197297
$ws.Cells.Item(14, 17).Value = 0.9687233804799998
$ws.Cells.Item(14, 18).Value = 8.71851042432
$ws.Cells.Item(14, 19).Value = 0.01091850334270543
$ws.Cells.Item(14, 20).Value = 0.01091850334270543

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Lama1"
$ws.Cells.Item(15, 3).Value = "Rpsa"
$ws.Cells.Item(15, 4).Value = "MuSCs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.006188666666666666
$ws.Cells.Item(15, 8).Value = 0.018566
$ws.Cells.Item(15, 9).Value = 0.07778257244722256
$ws.Cells.Item(15, 10).Value = 0.07778257244722254
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 389.2008666666666
$ws.Cells.Item(15, 14).Value = 1167.6026
$ws.Cells.Item(15, 15).Value = 0.3490212885850074
$ws.Cells.Item(15, 16).Value = 0.3490212885850074
$ws.Cells.Item(15, 17).Value = 2.408634430177777
$ws.Cells.Item(15, 18).Value = 21.6777098716
$ws.Cells.Item(15, 19).Value = 0.02714777366498631
$ws.Cells.Item(15, 20).Value = 0.0271477736649863

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Lama1"
$ws.Cells.Item(16, 3).Value = "Rpsa"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.006188666666666666
$ws.Cells.Item(16, 8).Value = 0.018566
$ws.Cells.Item(16, 9).Value = 0.07778257244722256
$ws.Cells.Item(16, 10).Value = 0.07778257244722254
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 169.4499613333333
$ws.Cells.Item(16, 14).Value = 508.349884
$ws.Cells.Item(16, 15).Value = 0.1519566088373896
$ws.Cells.Item(16, 16).Value = 0.1519566088373896
$ws.Cells.Item(16, 17).Value = 1.048669327371555
$ws.Cells.Item(16, 18).Value = 9.438023946344
$ws.Cells.Item(16, 19).Value = 0.01181957593572851
$ws.Cells.Item(16, 20).Value = 0.01181957593572851

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Lama1"
$ws.Cells.Item(17, 3).Value = "Rpsa"
$ws.Cells.Item(17, 4).Value = "FAPs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.006188666666666666
$ws.Cells.Item(17, 8).Value = 0.018566
$ws.Cells.Item(17, 9).Value = 0.07778257244722256
$ws.Cells.Item(17, 10).Value = 0.07778257244722254
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 169.0002543333333
$ws.Cells.Item(17, 14).Value = 507.000763
$ws.Cells.Item(17, 15).Value = 0.1515533278324679
$ws.Cells.Item(17, 16).Value = 0.1515533278324679
$ws.Cells.Item(17, 17).Value = 1.045886240650889
$ws.Cells.Item(17, 18).Value = 9.412976165858
$ws.Cells.Item(17, 19).Value = 0.01178820770174661
$ws.Cells.Item(17, 20).Value = 0.01178820770174661

# Row 18
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Lama1"
$ws.Cells.Item(18, 3).Value = "Rpsa"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.006188666666666666
$ws.Cells.Item(18, 8).Value = 0.018566
$ws.Cells.Item(18, 9).Value = 0.07778257244722256
$ws.Cells.Item(18, 10).Value = 0.07778257244722254
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 60.30985666666667
$ws.Cells.Item(18, 14).Value = 180.92957
$ws.Cells.Item(18, 15).Value = 0.05408370250677011
$ws.Cells.Item(18, 16).Value = 0.05408370250677011
$ws.Cells.Item(18, 17).Value = 0.3732375996244444
$ws.Cells.Item(18, 18).Value = 3.35913839662
$ws.Cells.Item(18, 19).Value = 0.004206769508446878
$ws.Cells.Item(18, 20).Value = 0.004206769508446877

# Row 19
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Lama1"
$ws.Cells.Item(19, 3).Value = "Rpsa"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.006188666666666666
$ws.Cells.Item(19, 8).Value = 0.018566
$ws.Cells.Item(19, 9).Value = 0.07778257244722256
$ws.Cells.Item(19, 10).Value = 0.07778257244722254
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 170.6279296666667
$ws.Cells.Item(19, 14).Value = 511.883789
$ws.Cells.Item(19, 15).Value = 0.1530129683186351
$ws.Cells.Item(19, 16).Value = 0.1530129683186351
$ws.Cells.Item(19, 17).Value = 1.055959380730444
$ws.Cells.Item(19, 18).Value = 9.503634426574
$ws.Cells.Item(19, 19).Value = 0.01190174229360881
$ws.Cells.Item(19, 20).Value = 0.01190174229360881
